$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row of data for 21 April 2020 (Excel serial date 43942)
$ws.Range("A38").Value = 43942
$ws.Range("A38").NumberFormat = "YYYY\-MM\-DD"

$ws.Range("B38").Value = -536
$ws.Range("D38").Value = -41
$ws.Range("F38").Value = 474
$ws.Range("G38").Value = 1366

# Update the active selection to mirror the original author's cursor position
$ws.Range("F39").Select()
